# Add the "associatedReferences" column to the "synonyms" sheet, between
# "acceptedNamecode" (col E) and "notes" (col F), mirroring the column that
# already exists in the "accepted names" sheet.

$wb = $excel.ActiveWorkbook
$synonyms = $wb.Worksheets.Item("synonyms")

# Insert a new column at F (pushes the existing "notes" column from F to G).
$synonyms.Columns.Item(6).Insert()

# New header cell + column width (~21.14 "characters", matching the
# associatedReferences column width used on the "accepted names" sheet).
$synonyms.Range("F1").Value = "associatedReferences"
$synonyms.Columns.Item(6).ColumnWidth = 20.43

# Make "synonyms" the active sheet / tab and move the in-sheet selection.
$synonyms.Activate()
$synonyms.Range("F3").Select()
